$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in column A (rows 2 and 3): 2025-12-30 -> 2025-12-31
$ws.Range("A2").Value = (Get-Date -Year 2025 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A3").Value = (Get-Date -Year 2025 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0).Date

# Update ticket_id values in column B (rows 2 and 3)
$ws.Range("B2").Value = 424267123
$ws.Range("B3").Value = 430947123

# Widen column B so the longer ticket_id values continue to fit
# (mirrors Excel's best-fit column-width recalculation)
$ws.Columns("B").ColumnWidth = 8.83

# Move the active selection to B4
$ws.Range("B4").Select()
